# Apply data updates for AfDD_2023_Annex_Table_Tab08.xlsx (Tab08 sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tab08")

# Updated source citation (World Bank data release date refresh)
$ws.Range("A104").Value = "Source: International Labour Organisation (retrieved 26/09/2022), United Nations Statistics Division National Accounts (Analysis of Main Aggregates dataset uploaded in January 2023), World Bank World Development Indicators (database and data releases from central banks, national statistical agencies, and World Bank country desks -10/10/2023)."

# Row 85 (Sao Tome and Principe) previously had no data ("..") for C/D/E; now populated
$ws.Range("C85").Value = 0.01786950565819
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0.42852544964292

# Refreshed indicator values (recalculated aggregates, rows 13-98)
$ws.Range("C13").Value = 4.9069218329992701
$ws.Range("D13").Value = 4.0702138176110401
$ws.Range("E13").Value = 12.3050871060241
$ws.Range("C23").Value = 7.1741656781693903
$ws.Range("D23").Value = 9.4489297352117898
$ws.Range("E23").Value = 22.6963769011004
$ws.Range("C38").Value = 0.52328773099134995
$ws.Range("D38").Value = 1.71818251345766
$ws.Range("E38").Value = 6.0555367833267901
$ws.Range("C45").Value = 7.7761845521380799
$ws.Range("D45").Value = 0.13235854750083001
$ws.Range("E45").Value = 11.058046646110199
$ws.Range("C61").Value = 4.3790654772602897
$ws.Range("D61").Value = 1.84605745804588
$ws.Range("E61").Value = 8.9620183256069197
$ws.Range("C62").Value = 5.1723222678423104
$ws.Range("D62").Value = 1.9879265659559
$ws.Range("E62").Value = 10.4893206119652
$ws.Range("C63").Value = 1.6827480315485801
$ws.Range("D63").Value = 0.91528348296096995
$ws.Range("E63").Value = 3.73243902084868
$ws.Range("C64").Value = 2.0216591903926999
$ws.Range("D64").Value = 3.4899234459460899
$ws.Range("E64").Value = 6.1407057827970402
$ws.Range("C65").Value = 1.23988388192919
$ws.Range("D65").Value = 0.89249697280849005
$ws.Range("E65").Value = 3.6593282005323702
$ws.Range("C66").Value = 1.86223096331427
$ws.Range("D66").Value = 0.97045175821680996
$ws.Range("E66").Value = 4.0799595609047499
$ws.Range("C67").Value = 4.4247636842766003
$ws.Range("D67").Value = 2.22438644588722
$ws.Range("E67").Value = 9.9376539568048301
$ws.Range("C68").Value = 5.2426171432869699
$ws.Range("D68").Value = 1.19191826001304
$ws.Range("E68").Value = 8.7899519060320408
$ws.Range("C69").Value = 0.10250538868568
$ws.Range("D69").Value = 5.5590828634302296
$ws.Range("E69").Value = 9.8173481651564298
$ws.Range("C70").Value = 14.4979564905323
$ws.Range("D70").Value = 5.4853271854221397
$ws.Range("E70").Value = 24.34468503387
$ws.Range("C71").Value = 4.3790654772602897
$ws.Range("D71").Value = 1.84605745804588
$ws.Range("E71").Value = 8.9620183256069108
$ws.Range("C72").Value = 0.68770580587665997
$ws.Range("D72").Value = 1.47182955003094
$ws.Range("E72").Value = 6.2506664460432999
$ws.Range("C73").Value = 3.8721477825519699
$ws.Range("D73").Value = 5.5243743037425004
$ws.Range("E73").Value = 13.049740790379101
$ws.Range("C74").Value = 13.7379424265708
$ws.Range("D74").Value = 0.29708269570092
$ws.Range("E74").Value = 18.4243267702408
$ws.Range("C75").Value = 22.856002386360199
$ws.Range("D75").Value = 0.021920976814660001
$ws.Range("E75").Value = 26.996474700211401
$ws.Range("C76").Value = 0.69267360955753998
$ws.Range("D76").Value = 1.01099382470901
$ws.Range("E76").Value = 3.92272189496367
$ws.Range("C77").Value = 2.2654977569793102
$ws.Range("D77").Value = 4.6270117009172402
$ws.Range("E77").Value = 7.6916214033411503
$ws.Range("C78").Value = 0.046253008960359998
$ws.Range("D78").Value = 0.081447887142499995
$ws.Range("E78").Value = 0.26655302563244998
$ws.Range("C79").Value = 0.49783272968526998
$ws.Range("D79").Value = 0.53996194483394
$ws.Range("E79").Value = 1.4236368448443599
$ws.Range("C80").Value = 23.486153272542801
$ws.Range("D80").Value = 0.25061551747157002
$ws.Range("E80").Value = 29.4976430091662
$ws.Range("C81").Value = 15.325050579753899
$ws.Range("D81").Value = 2.17457612076814
$ws.Range("E81").Value = 23.0973320013021
$ws.Range("C82").Value = 2.1938149062069101
$ws.Range("D82").Value = 2.2704777015447801
$ws.Range("E82").Value = 7.39786322230872
$ws.Range("C83").Value = 0.50900241847322003
$ws.Range("D83").Value = 0.80694218623084002
$ws.Range("E83").Value = 2.0664103682903101
$ws.Range("C84").Value = 0.94613430401349996
$ws.Range("D84").Value = 6.7224842385492103
$ws.Range("E84").Value = 13.691400914286399
$ws.Range("C86").Value = 5.5452177393205702
$ws.Range("D86").Value = 0.57199181020681
$ws.Range("E86").Value = 8.7954510485540496
$ws.Range("C87").Value = 1.5562444008169101
$ws.Range("D87").Value = 1.6393815568732
$ws.Range("E87").Value = 5.41173501015634
$ws.Range("C88").Value = 8.1528443793115404
$ws.Range("D88").Value = 2.9130944240727299
$ws.Range("E88").Value = 14.269465169760499
$ws.Range("C89").Value = 2.1933182537992
$ws.Range("D89").Value = 1.19026350487406
$ws.Range("E89").Value = 4.8453078661057098
$ws.Range("C90").Value = 1.3514283522608901
$ws.Range("D90").Value = 0.47052501728516
$ws.Range("E90").Value = 2.35271319025113
$ws.Range("C91").Value = 4.1113557020457296
$ws.Range("D91").Value = 5.3019900459475897
$ws.Range("E91").Value = 14.2154423124449
$ws.Range("C92").Value = 0.05281138161915
$ws.Range("D92").Value = 0.400221783045
$ws.Range("E92").Value = 1.9977493595271401
$ws.Range("E93").Value = 2.92383933181728
$ws.Range("C94").Value = 0.65467996405868001
$ws.Range("D94").Value = 1.0491634833535399
$ws.Range("E94").Value = 2.53762386306123
$ws.Range("C95").Value = 0.56503505265978005
$ws.Range("D95").Value = 4.6186007331084298
$ws.Range("E95").Value = 10.4095265513882
$ws.Range("C96").Value = 6.9403401826854898
$ws.Range("D96").Value = 6.1452953171836002
$ws.Range("E96").Value = 17.267108020076702
$ws.Range("C97").Value = 6.9176896816407396
$ws.Range("D97").Value = 1.69046656587552
$ws.Range("E97").Value = 11.8892971012302
$ws.Range("C98").Value = 5.7877018880995701
$ws.Range("D98").Value = 1.44901653412706
$ws.Range("E98").Value = 10.9509752491716
